$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: Target cluster changes from "Resolving-Mac" to "ECs", and the
# TPM-derived numeric columns (G, H, M, N, O, P, Q, R, S, T) are refreshed
# with newly computed values.
$ws.Range("D2").Value = "ECs"
$ws.Range("G2").Value = 1.030436666666667
$ws.Range("H2").Value = 3.09131
$ws.Range("M2").Value = 0.016584
$ws.Range("N2").Value = 0.049752
$ws.Range("O2").Value = 0.05771573215408494
$ws.Range("P2").Value = 0.05771573215408494
$ws.Range("Q2").Value = 0.01708876168
$ws.Range("R2").Value = 0.15379885512
$ws.Range("S2").Value = 0.05771573215408494
$ws.Range("T2").Value = 0.05771573215408494

# New row 3: same sending cluster / ligand / receptor, target cluster "MuSCs".
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vip"
$ws.Range("C3").Value = "Sctr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.030436666666667
$ws.Range("H3").Value = 3.09131
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.039802
$ws.Range("N3").Value = 0.119406
$ws.Range("O3").Value = 0.1385191492521038
$ws.Range("P3").Value = 0.1385191492521038
$ws.Range("Q3").Value = 0.04101344020666666
$ws.Range("R3").Value = 0.36912096186
$ws.Range("S3").Value = 0.1385191492521038
$ws.Range("T3").Value = 0.1385191492521038

# New row 4: same sending cluster / ligand / receptor, target cluster
# "Resolving-Mac".
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vip"
$ws.Range("C4").Value = "Sctr"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.030436666666667
$ws.Range("H4").Value = 3.09131
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2309533333333333
$ws.Range("N4").Value = 0.69286
$ws.Range("O4").Value = 0.8037651185938113
$ws.Range("P4").Value = 0.8037651185938113
$ws.Range("Q4").Value = 0.2379827829555556
$ws.Range("R4").Value = 2.1418450466
$ws.Range("S4").Value = 0.8037651185938113
$ws.Range("T4").Value = 0.8037651185938113
